$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2, column N ("Comments") was left blank - fill in the missing comment
# value, matching the formatting already used by the header cell N1.
$ws.Cells.Item(2, 14).Value = "Selenium Class"
$ws.Range("N2").Font.ThemeColor = 1
